$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2008.2683
$ws.Range("I15").Value = 2008.2683
$ws.Range("K15").Value = 6024.8049
$ws.Range("M15").Value = -5855.8049
$ws.Range("H19").Value = 563.9167
$ws.Range("I19").Value = 499.25
$ws.Range("K19").Value = 499.25
$ws.Range("M19").Value = -324.25
$ws.Range("H69").Value = 9334.25
$ws.Range("I69").Value = 7668
$ws.Range("J69").Value = 14333
$ws.Range("K69").Value = 23004
$ws.Range("L69").Value = 42999
$ws.Range("M69").Value = -22130
$ws.Range("N69").Value = -44747
$ws.Range("H70").Value = 1565.1666
$ws.Range("J70").Value = 1619.4286
$ws.Range("L70").Value = 4858.2858
$ws.Range("N70").Value = -5398.2858
$ws.Range("H72").Value = 9334.25
$ws.Range("I72").Value = 7668
$ws.Range("J72").Value = 14333
$ws.Range("K72").Value = 69012
$ws.Range("L72").Value = 128997
$ws.Range("M72").Value = -64644
$ws.Range("N72").Value = -137733
$ws.Range("H73").Value = 1565.1666
$ws.Range("J73").Value = 1619.4286
$ws.Range("L73").Value = 4858.2858
$ws.Range("N73").Value = -6730.2858
$ws.Range("H112").Value = 2216.2222
$ws.Range("J112").Value = 2216.2222
$ws.Range("L112").Value = 6648.6666
$ws.Range("N112").Value = -8864.6666
$ws.Range("H116").Value = 6458.7
$ws.Range("J116").Value = 6398.6665
$ws.Range("L116").Value = 6398.6665
$ws.Range("N116").Value = -13282.6665
$ws.Range("H137").Value = 20835202
$ws.Range("I137").Value = 22729168
$ws.Range("J137").Value = 1583.5
$ws.Range("K137").Value = 68187504
$ws.Range("L137").Value = 4750.5
$ws.Range("M137").Value = -68184954
$ws.Range("N137").Value = -9850.5
$ws.Range("H138").Value = 1888.569
$ws.Range("J138").Value = 2010.2979
$ws.Range("L138").Value = 6030.893700000001
$ws.Range("N138").Value = -16310.8937

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2468.0356
$ws.Range("I2").Value = 2393.4285
$ws.Range("J2").Value = 2542.6428
$ws.Range("K2").Value = 2393.4285
$ws.Range("L2").Value = 2542.6428
$ws.Range("M2").Value = -2280.4285
$ws.Range("N2").Value = -2768.6428
$ws.Range("H32").Value = 3493.652
$ws.Range("I32").Value = 1693.6
$ws.Range("J32").Value = 9221.091
$ws.Range("K32").Value = 1693.6
$ws.Range("L32").Value = 9221.091
$ws.Range("M32").Value = -1406.6
$ws.Range("N32").Value = -9795.091
$ws.Range("H45").Value = 3267.5833
$ws.Range("J45").Value = 3151.75
$ws.Range("L45").Value = 3151.75
$ws.Range("N45").Value = -3905.75
$ws.Range("H74").Value = 2520.8333
$ws.Range("I74").Value = 2080.1667
$ws.Range("K74").Value = 2080.1667
$ws.Range("M74").Value = -1206.1667
$ws.Range("H77").Value = 2520.8333
$ws.Range("I77").Value = 2080.1667
$ws.Range("K77").Value = 10400.8335
$ws.Range("M77").Value = -6032.833500000001
$ws.Range("H102").Value = 61202.8
$ws.Range("I102").Value = 12753.625
$ws.Range("J102").Value = 254999.5
$ws.Range("K102").Value = 12753.625
$ws.Range("L102").Value = 254999.5
$ws.Range("M102").Value = -11131.625
$ws.Range("N102").Value = -258243.5
$ws.Range("H110").Value = 3926.3
$ws.Range("I110").Value = 3909.375
$ws.Range("J110").Value = 3994
$ws.Range("K110").Value = 3909.375
$ws.Range("L110").Value = 3994
$ws.Range("M110").Value = -1864.375
$ws.Range("N110").Value = -8084
$ws.Range("H116").Value = 2468.0356
$ws.Range("I116").Value = 2393.4285
$ws.Range("J116").Value = 2542.6428
$ws.Range("K116").Value = 2393.4285
$ws.Range("L116").Value = 2542.6428
$ws.Range("M116").Value = -99.42849999999999
$ws.Range("N116").Value = -7130.6428
$ws.Range("H132").Value = 9052.666999999999
$ws.Range("I132").Value = 6496.364
$ws.Range("K132").Value = 19489.092
$ws.Range("M132").Value = -16959.092

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2468.0356
$ws.Range("I3").Value = 2393.4285
$ws.Range("J3").Value = 2542.6428
$ws.Range("K3").Value = 2393.4285
$ws.Range("L3").Value = 2542.6428
$ws.Range("M3").Value = -2279.4285
$ws.Range("N3").Value = -2770.6428
$ws.Range("H86").Value = 3067.6875
$ws.Range("I86").Value = 2814.875
$ws.Range("J86").Value = 3320.5
$ws.Range("K86").Value = 2814.875
$ws.Range("L86").Value = 3320.5
$ws.Range("M86").Value = -1691.875
$ws.Range("N86").Value = -5566.5
$ws.Range("H89").Value = 3067.6875
$ws.Range("I89").Value = 2814.875
$ws.Range("J89").Value = 3320.5
$ws.Range("K89").Value = 14074.375
$ws.Range("L89").Value = 16602.5
$ws.Range("M89").Value = -8458.375
$ws.Range("N89").Value = -27834.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 280.82758
$ws.Range("I22").Value = 264.125
$ws.Range("J22").Value = 301.3846
$ws.Range("K22").Value = 264.125
$ws.Range("L22").Value = 301.3846
$ws.Range("M22").Value = 85.875
$ws.Range("N22").Value = -1001.3846
$ws.Range("H59").Value = 35000
$ws.Range("I59").Value = 10000
$ws.Range("K59").Value = 10000
$ws.Range("M59").Value = -8855
$ws.Range("H99").Value = 82053800
$ws.Range("I99").Value = 133334460
$ws.Range("J99").Value = 66669596
$ws.Range("K99").Value = 133334460
$ws.Range("L99").Value = 66669596
$ws.Range("M99").Value = -133332962
$ws.Range("N99").Value = -66672592
$ws.Range("H122").Value = 2307.25
$ws.Range("I122").Value = 1996.25
$ws.Range("K122").Value = 5988.75
$ws.Range("M122").Value = -3538.75
$ws.Range("H126").Value = 82053800
$ws.Range("I126").Value = 133334460
$ws.Range("J126").Value = 66669596
$ws.Range("K126").Value = 400003380
$ws.Range("L126").Value = 200008788
$ws.Range("M126").Value = -400000910
$ws.Range("N126").Value = -200013728

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 668.1
$ws.Range("I8").Value = 668.1
$ws.Range("K8").Value = 2004.3
$ws.Range("M8").Value = -1865.3
$ws.Range("H46").Value = 88626250
$ws.Range("J46").Value = 95240340
$ws.Range("L46").Value = 285721020
$ws.Range("N46").Value = -285721202
$ws.Range("H58").Value = 800
$ws.Range("I58").Value = 800
$ws.Range("J58").Value = 800
$ws.Range("K58").Value = 2400
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -2272
$ws.Range("N58").Value = -2656
$ws.Range("H99").Value = 3366.3333
$ws.Range("I99").Value = 2071.6
$ws.Range("K99").Value = 6214.799999999999
$ws.Range("M99").Value = -3968.799999999999
$ws.Range("H108").Value = 1431
$ws.Range("I108").Value = 1431
$ws.Range("K108").Value = 4293
$ws.Range("M108").Value = -1413
$ws.Range("H115").Value = 236684.67
$ws.Range("I115").Value = 352527.5
$ws.Range("J115").Value = 4999
$ws.Range("K115").Value = 1057582.5
$ws.Range("L115").Value = 14997
$ws.Range("M115").Value = -1056407.5
$ws.Range("N115").Value = -17347
$ws.Range("H120").Value = 12749.833
$ws.Range("I120").Value = 9499.666999999999
$ws.Range("K120").Value = 28499.001
$ws.Range("M120").Value = -23661.001
$ws.Range("H132").Value = 775.75
$ws.Range("I132").Value = 775.75
$ws.Range("K132").Value = 6981.75
$ws.Range("M132").Value = -4451.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 7349.5
$ws.Range("I6").Value = 3500
$ws.Range("J6").Value = 8632.666999999999
$ws.Range("K6").Value = 3500
$ws.Range("L6").Value = 8632.666999999999
$ws.Range("N6").Value = -8858.666999999999
$ws.Range("M6").Value = -3387
$ws.Range("H16").Value = 7349.5
$ws.Range("I16").Value = 3500
$ws.Range("J16").Value = 8632.666999999999
$ws.Range("K16").Value = 3500
$ws.Range("L16").Value = 8632.666999999999
$ws.Range("N16").Value = -9132.666999999999
$ws.Range("M16").Value = -3250
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H70").Value = 24953.5
$ws.Range("I70").Value = 36520.65
$ws.Range("K70").Value = 36520.65
$ws.Range("M70").Value = -36250.65
$ws.Range("H73").Value = 24953.5
$ws.Range("I73").Value = 36520.65
$ws.Range("K73").Value = 36520.65
$ws.Range("M73").Value = -35584.65

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 19999
$ws.Range("I12").Value = 19999
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 19999
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -19829
$ws.Range("N12").ClearContents()
$ws.Range("H22").Value = 6932.3335
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 8898.5
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 8898.5
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -9488.5
$ws.Range("H27").Value = 6932.3335
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 8898.5
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 8898.5
$ws.Range("M27").Value = -2893
$ws.Range("N27").Value = -9112.5
$ws.Range("H55").Value = 336.33334
$ws.Range("J55").Value = 232.5
$ws.Range("L55").Value = 232.5
$ws.Range("N55").Value = -578.5
$ws.Range("H132").Value = 2553.55
$ws.Range("I132").Value = 2323.2812
$ws.Range("K132").Value = 6969.8436
$ws.Range("M132").Value = -4439.8436

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13212.429
$ws.Range("J41").Value = 12811
$ws.Range("L41").Value = 12811
$ws.Range("N41").Value = -13591
$ws.Range("H122").Value = 20458484
$ws.Range("I122").Value = 19234136
$ws.Range("K122").Value = 57702408
$ws.Range("M122").Value = -57699958
$ws.Range("H123").Value = 39869.5
$ws.Range("J123").Value = 57250
$ws.Range("L123").Value = 57250
$ws.Range("N123").Value = -67050
$ws.Range("H136").Value = 3633.8333
$ws.Range("I136").Value = 2674.5715
$ws.Range("K136").Value = 8023.7145
$ws.Range("M136").Value = -5473.7145

Write-Host "Applied all changes"
